# Diary workbook update: add the "2 marras" entry as row 22.
#
# Shared-string insertion order matters for this engine (new strings are
# appended to sharedStrings.xml in the order the cell values are first
# assigned), so columns are written in the same order the target sheet's
# string table expects: A, B, C, E, D, F, G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A22 - date
$ws.Range("A22").Value = "2 marras"

# B22 - time range (time-of-day format, wrapped, like the other "Kello" cells)
$ws.Range("B22").Value = "18.45-22.15"
$ws.Range("B22").NumberFormat = "h:mm"
$ws.Range("B22").WrapText = $true

# C22 - learning content (reuses the existing "Kangassimulaatio" string)
$ws.Range("C22").Value = "Kangassimulaatio"
$ws.Range("C22").WrapText = $true

# E22 before D22 so new shared strings land in the same order as the target file
$ws.Range("E22").Value = "Hyvin taas opittu c++ kummallisuuksia ja linkkerinkin toimintaa. Ihan hyvä meno tuon demonkin kanssa, huomenna se toimi (:"
$ws.Range("E22").WrapText = $true

$ws.Range("D22").Value = "Juujuu, tällaista se on kun ei tiedä mitä ei tiedä ja tutoriaalikoodia modernisoiden kompastellaan."
$ws.Range("D22").WrapText = $true

$ws.Range("F22").Value = "Puolivälipaniikki?"
$ws.Range("F22").WrapText = $true

# G22 - hours logged for the session
$ws.Range("G22").Value = 3.5

# Row 22 ends up the same height as other wrapped-text rows in the log (58pt).
$ws.Rows.Item(22).RowHeight = 58

# Move selection to the newly added row, matching where the author ended up.
$ws.Range("F22").Select()
